$wb = $excel.ActiveWorkbook

# The same set of updates needs to be applied to both the "展览" and
# "全部类型" worksheets, which contain duplicate data tables.
$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F ("想去人数")
$updates = @{
    3  = 7884
    5  = 193
    9  = 128
    10 = 175
    12 = 714
    14 = 1875
    16 = 55
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
